# Daily attendance processing - reverse the order of the comma-separated
# "Recorded By" entries (column G) for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in the sheet.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $reversedParts = $parts[($parts.Count - 1)..0]
            $cell.Value2 = [string]::Join(", ", $reversedParts)
        }
    }
}
